$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.823.53"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "3.850.86"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'598.09"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'166.19"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "3.847.63"
$ws.Range("E7").Value = "  -1.63%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "'6.31"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "'36.78"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "4.498.28"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").Value = "3.836.39"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "67.871.56"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "'18.10"
$ws.Range("E18").Value = "  +6.58%  "
$ws.Range("D19").Value = "'7.38"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("D22").Value = "'465.41"
$ws.Range("E22").Value = "  -3.71%  "
$ws.Range("D23").Value = "'0.727"
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("E24").Value = "  -4.46%  "
$ws.Range("D25").Value = "'83.24"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").Value = "'12.10"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'9.98"
$ws.Range("D30").Value = "'2.94"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "4.002.93"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").Value = "'2.30"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").Value = "'30.94"
$ws.Range("E34").Value = "  -3.05%  "
$ws.Range("D35").Value = "3.828.57"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'1.01"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("D39").Value = "'3.28"
$ws.Range("E39").Value = "  +8.64%  "
$ws.Range("D40").Value = "'5.88"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'0.311"
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").Value = "'427.12"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'47.20"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").Value = "'8.52"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").Value = "'143.76"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").Value = "'0.000270"
$ws.Range("E49").Value = "  +6.33%  "
$ws.Range("D50").Value = "'25.73"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").Value = "'39.24"
$ws.Range("E51").Value = "  +0.44%  "
